# Auto-generated edit script applying updated "想去人数" (F) and "最低票价" (G) values
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 1862
$ws1.Range("G4").Value = 68.8
$ws1.Range("F7").Value = 83
$ws1.Range("F8").Value = 88
$ws1.Range("F9").Value = 237
$ws1.Range("F10").Value = 171
$ws1.Range("F11").Value = 1089
$ws1.Range("F12").Value = 360
$ws1.Range("F13").Value = 86
$ws1.Range("F14").Value = 58
$ws1.Range("F15").Value = 108
$ws1.Range("F17").Value = 233
$ws1.Range("F20").Value = 1200
$ws1.Range("F21").Value = 464
$ws1.Range("F25").Value = 542
$ws1.Range("F27").Value = 52
$ws1.Range("F28").Value = 1916
$ws1.Range("F29").Value = 2375
$ws1.Range("F30").Value = 1155
$ws1.Range("F32").Value = 82
$ws1.Range("F34").Value = 325
$ws1.Range("F35").Value = 714
$ws1.Range("F36").Value = 656
$ws1.Range("F37").Value = 89
$ws1.Range("F39").Value = 711
$ws1.Range("F40").Value = 118
$ws1.Range("F41").Value = 519
$ws1.Range("F42").Value = 577
$ws1.Range("F43").Value = 263
$ws1.Range("F44").Value = 183

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F8").Value = 164
$ws2.Range("F15").Value = 210
$ws2.Range("F20").Value = 5
$ws2.Range("F21").Value = 4

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1862
$ws4.Range("G4").Value = 68.8
$ws4.Range("F7").Value = 83
$ws4.Range("F10").Value = 89
$ws4.Range("F11").Value = 237
$ws4.Range("F12").Value = 171
$ws4.Range("F14").Value = 1089
$ws4.Range("F15").Value = 360
$ws4.Range("F16").Value = 86
$ws4.Range("F17").Value = 58
$ws4.Range("F19").Value = 233
$ws4.Range("F22").Value = 1200
$ws4.Range("F23").Value = 464
$ws4.Range("F27").Value = 52
$ws4.Range("F28").Value = 1916
$ws4.Range("F29").Value = 2375
$ws4.Range("F31").Value = 1155
$ws4.Range("F35").Value = 82
$ws4.Range("F37").Value = 325
$ws4.Range("F38").Value = 5
$ws4.Range("F39").Value = 4
$ws4.Range("F40").Value = 714
$ws4.Range("F41").Value = 656
$ws4.Range("F42").Value = 711
$ws4.Range("F43").Value = 118
$ws4.Range("F44").Value = 519
$ws4.Range("F45").Value = 577
$ws4.Range("F46").Value = 263
$ws4.Range("F48").Value = 183

